$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 53 values
$ws.Range("C53").Value = 12.5
$ws.Range("F53").Value = 26.5

# Add new row 54 with the latest quarter data
# Force the period label to stay text (it looks like a date and would
# otherwise be auto-converted to a date serial by Excel's type inference)
$ws.Range("A54").NumberFormat = "@"
$ws.Range("A54").Value = "01-04-2021"
$ws.Range("A54").Style = "Normal"

$ws.Range("B54").Value = 3.4
$ws.Range("C54").Value = 16
$ws.Range("D54").Value = 15.5
$ws.Range("E54").Value = 0.4
$ws.Range("F54").Value = 24.3
$ws.Range("G54").Value = 8
$ws.Range("H54").Value = 0.1
$ws.Range("I54").Value = 4.5
$ws.Range("J54").Value = 27.9
